$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")
$ws.Activate() | Out-Null
$ws.Range("J8").Select() | Out-Null
